$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "process" column (column Y) with header and a constant
# value for every data row, matching the new shared-strings entries.
$ws.Range("Y1").Value = "process"
$ws.Range("Y2:Y35").Value = "pp->Z/gamma*->l+ l-"

# Reflect the saved selection/view state from the edit: column Y is
# selected (Y1:Y35) with Y1 as the active cell.
$ws.Range("Y1:Y35").Select()
